$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "As a data analyst in an educational institution, you've been given a large Excel spreadsheet of student names along with their respective test scores, shown below. The Dean has asked you to identify the student with the highest test score.Which tool should you use to achieve this most efficiently?",
        "ques_type": 2,
        "options": [
            "Filter",
            "Sort",
            "Find",
            "Search"
        ],
        "score": "Sort"
    },
    {
        "title": "As a project manager handling the data for a technical team's project in Excel, you've been asked to ensure all numerical data displays 10 decimal points. Currently, only six decimal points are being displayed.What should you adjust to achieve this?",
        "ques_type": 2,
        "options": [
            "The Number format",
            "The Percentage format",
            "The Fraction format",
            "The Special format"
        ],
        "score": "The Number format"
    },
    {
        "title": "As a financial analyst at your company, you're scrutinizing an Excel spreadsheet detailing contractor fees. You've noticed that the fees are currently shown as plain numbers, but they should be shown in dollar amounts.Which format type should you select for these cells?",
        "ques_type": 2,
        "options": [
            "Accounting",
            "Number",
            "Special",
            "Scientific"
        ],
        "score": "Accounting"
    },
    {
        "title": "As an administrative assistant, you're managing an Excel spreadsheet that tracks staff training, shown below. You need to rename the tab labeled Sheet1 to Training Dates.Which action should you take to achieve this most efficiently?",
        "ques_type": 2,
        "options": [
            "Right-click on the Sheet1 tab &gt click Rename",
            "Click File &gt select Save As",
            "Press Shift (PC)/Function (Apple) + F11",
            "Select Sheet Options in the Page Layout tab"
        ],
        "score": "Right-click on the Sheet1 tab &gt click Rename"
    }
]
'@

$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
$ws.Rows.Item(1).AutoFit()
